$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: split the run that begins the "(2) Submit copies of D&Fs ..."
#           paragraph into two runs: "(2)" and " Submit copies of D&Fs ...".
#           A temporary bookmark is inserted/removed exactly at the split
#           point; Word's run-coalescing only merges runs that are fully
#           identical, so round-tripping a bookmark through that location
#           forces the paragraph's run list to break there without
#           introducing any stray run-level formatting.
# ---------------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute( `
    "(2) Submit copies of D&Fs executed pursuant to DFARS 245.102(4)(ii)(C)(1) to ", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $prefix = "(2)"
    $splitPoint = $rng.Start + $prefix.Length
    $boundary = $d.Range($splitPoint, $splitPoint)
    $bm = $d.Bookmarks.Add("zzTempSplitMark", $boundary)
    $d.Bookmarks("zzTempSplitMark").Delete()
}

# ---------------------------------------------------------------------------
# Change 2: the paragraph holding the eight leading spaces (immediately
#           followed by the tab + "(2) Provide retention and redistribution
#           ..." run) needs the List 2 paragraph style applied.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Provide retention and redistribution requirements of the owning Commands*") {
        $p.Range.Style = $d.Styles("List 2")
    }
}
